# Applies the edits described by the commit:
#  - workbook window size tweak (best effort; engine may not persist this)
#  - rewrite the "总结：" summary cell (A80) with the fuller wrap-up text
#  - mark the last five task rows (C74:C78) as "完成" (Completed)
#  - row 76 gets a slightly shorter custom height (37 instead of 40)
#  - move the active selection from A80:D80 down to F80, scrolled to A70

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- window geometry (best effort) ---------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.Width = 19095
    $win.Height = 12210
} catch {
    # Older/limited hosts may not expose window sizing; ignore.
}

# --- update the final summary note ----------------------------------------
$ws.Range("A80").Value = "总结：虽然勉勉强强把原型做好了，但接下来要进行开发，才是重点，要加油。"

# --- mark the "完成情况" column as done for rows 74-78 ---------------------
$ws.Range("C74").Value = "完成"
$ws.Range("C75").Value = "完成"
$ws.Range("C76").Value = "完成"
$ws.Range("C77").Value = "完成"
$ws.Range("C78").Value = "完成"

# --- tighten row 76's custom height from 40 to 37 --------------------------
$ws.Rows.Item(76).RowHeight = 37

# --- scroll the view up slightly and move the selection to F80 -------------
try {
    $excel.ActiveWindow.ScrollRow = 70
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Ignore if scrolling isn't supported.
}
$ws.Range("F80").Select() | Out-Null
